# Daily attendance processing - 2025-10-26 13:44:31
# Rotate the "Recorded By" (column G) name/email list for each data row:
# the first name in the comma-separated list moves to the end of the list.
# e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # -4162 = xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if (-not ($val -is [string])) { continue }
    if ($val.Length -eq 0) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ","
    for ($i = 0; $i -lt $parts.Length; $i++) {
        $parts[$i] = $parts[$i].Trim()
    }

    $rotated = @($parts[1..($parts.Length - 1)]) + @($parts[0])
    $cell.Value = [string]::Join(", ", $rotated)
}
